$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.118.41'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.759.52'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '576.39'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.10'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  -3.11%  '
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.87'
$ws.Range('E10').Value = '  -13.78%  '
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('D13').Value = '3.248.94'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.99'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').Value = '63.735.81'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('E16').Value = '  -5.35%  '
$ws.Range('D17').Value = '2.765.09'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.12'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.82'
$ws.Range('E19').Value = '  -4.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '359.15'
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('E21').Value = '  -6.06%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.528'
$ws.Range('E23').Value = '  -8.49%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.07'
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('E25').Value = '  -3.33%  '
$ws.Range('E26').Value = '  -3.08%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '0.0₃0905'
$ws.Range('E28').Value = '  -6.72%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.38'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('E30').Value = '  -4.73%  '
$ws.Range('E31').Value = '  +6.14%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '169.49'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.94'
$ws.Range('E33').Value = '  -5.56%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '20.21'
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('E38').Value = '  -2.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '348.89'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.32'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  -2.44%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.13'
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '21.56'
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '22.06'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0589'
$ws.Range('E45').Value = '  -3.63%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '137.73'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.629'
$ws.Range('E47').Value = '  -3.72%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0253'
$ws.Range('E48').Value = '  -3.54%  '
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.04'
$ws.Range('E51').Value = '  +0.04%  '
